# Apply updates to the "mean", "stdev", and "summary" worksheets
# per the commit "Update basaltic glass analysis - major, wd scans and D2872"

$wb = $excel.ActiveWorkbook

# ---- Sheet "mean" ----
$wsMean = $wb.Worksheets.Item("mean")
$wsMean.Range("J3").Value = 5.169
$wsMean.Range("K3").Value = 0.182
$wsMean.Range("J4").Value = 5.132
$wsMean.Range("K4").Value = 0.219

# ---- Sheet "stdev" ----
$wsStdev = $wb.Worksheets.Item("stdev")
$wsStdev.Range("J3").Value = 0.025
$wsStdev.Range("K3").Value = 0.025
$wsStdev.Range("J4").Value = 0.031
$wsStdev.Range("K4").Value = 0.031

# ---- Sheet "summary" ----
$wsSummary = $wb.Worksheets.Item("summary")
$wsSummary.Range("L3").Value = 5.169
$wsSummary.Range("M3").Value = 0.025
$wsSummary.Range("T3").Value = 0.182
$wsSummary.Range("U3").Value = 0.025
$wsSummary.Range("L4").Value = 5.132
$wsSummary.Range("M4").Value = 0.031
$wsSummary.Range("T4").Value = 0.219
$wsSummary.Range("U4").Value = 0.031
